$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / label text updates ---
$ws.Range("A8").Value = "Volume 30   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/26/2022  Through  1/1/2023"
$ws.Range("M12").Value = "13 Year (2010)"
$ws.Range("N12").Value = "30 Year (1993)"
$ws.Range("K35").Value = "'22 vs '01"
$ws.Range("L35").Value = "'22 vs '98"
$ws.Range("M35").Value = "'22 vs '93"
$ws.Range("N35").Value = "'22 vs '90"

# --- Numeric cell value updates (style unchanged) ---
$ws.Range("C13").Value = 2023
$ws.Range("D13").Value = 2022
$ws.Range("F13").Value = 2023
$ws.Range("G13").Value = 2022
$ws.Range("I13").Value = 2023
$ws.Range("J13").Value = 2022
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("N15").Value = -100
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = -100
$ws.Range("N16").Value = -100
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 9.090909090909
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 100
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = 100
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 15
$ws.Range("E18").Value = -73.333333333333
$ws.Range("F18").Value = 35
$ws.Range("G18").Value = 44
$ws.Range("H18").Value = -20.454545454545
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = -100
$ws.Range("L18").Value = -100
$ws.Range("M18").Value = -100
$ws.Range("N18").Value = -100
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = 16.666666666666
$ws.Range("F19").Value = 110
$ws.Range("G19").Value = 107
$ws.Range("H19").Value = 2.803738317757
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -87.5
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 10
$ws.Range("N20").Value = -100
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -22.5
$ws.Range("F21").Value = 178
$ws.Range("G21").Value = 186
$ws.Range("H21").Value = -4.301075268817
$ws.Range("I21").Value = 3
$ws.Range("J21").Value = 6
$ws.Range("K21").Value = -50
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = -90.909090909090
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -66.666666666666
$ws.Range("C24").Value = 66
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = 29.411764705882
$ws.Range("F24").Value = 264
$ws.Range("G24").Value = 224
$ws.Range("H24").Value = 17.857142857142
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 3
$ws.Range("K24").Value = -66.666666666666
$ws.Range("L24").Value = -50
$ws.Range("M24").Value = 0
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 7
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = 30.769230769230
$ws.Range("M25").Value = -100
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = -100
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("J35").Value = 2022
$ws.Range("J37").Value = 12
$ws.Range("K37").Value = 33.333333333333
$ws.Range("L37").Value = -29.411764705882
$ws.Range("M37").Value = -60
$ws.Range("N37").Value = -55.555555555555
$ws.Range("J38").Value = 229
$ws.Range("K38").Value = -36.211699164345
$ws.Range("L38").Value = -50.108932461873
$ws.Range("M38").Value = -83.381712626995
$ws.Range("N38").Value = -87.972689075630
$ws.Range("J39").Value = 174
$ws.Range("K39").Value = 37.007874015748
$ws.Range("L39").Value = -5.945945945945
$ws.Range("M39").Value = -41.016949152542
$ws.Range("N39").Value = -44.230769230769
$ws.Range("J40").Value = 275
$ws.Range("K40").Value = -60.770328102710
$ws.Range("L40").Value = -73.249027237354
$ws.Range("M40").Value = -90.891023517721
$ws.Range("N40").Value = -92.557510148849
$ws.Range("J41").Value = 1719
$ws.Range("K41").Value = -22.146739130434
$ws.Range("L41").Value = -20.673742501153
$ws.Range("M41").Value = -54.667721518987
$ws.Range("N41").Value = -68.201997780244
$ws.Range("J42").Value = 191
$ws.Range("K42").Value = -40.3125
$ws.Range("L42").Value = -80.842527582748
$ws.Range("M42").Value = -94.474978304888
$ws.Range("N42").Value = -94.278010784901
$ws.Range("J43").Value = 2602
$ws.Range("K43").Value = -30.147651006711
$ws.Range("L43").Value = -46.449886807985
$ws.Range("M43").Value = -78.273213092852
$ws.Range("N43").Value = -82.288475937648

# --- C23: was text "0" (style 14), becomes numeric 1 (style 15) ---
$ws.Range("C23").Value = 1
$ws.Range("F23").Copy()
$ws.Range("C23").PasteSpecial(-4122)

# --- Cells converting from numeric to text placeholders ("0" / "***.*") ---
# Force text storage via Text number format, then restore style 14 via PasteSpecial formats.
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "0"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "0"
$ws.Range("K14").NumberFormat = "@"
$ws.Range("K14").Value = "***.*"
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value = "***.*"
$ws.Range("M14").NumberFormat = "@"
$ws.Range("M14").Value = "***.*"
$ws.Range("N14").NumberFormat = "@"
$ws.Range("N14").Value = "***.*"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "0"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "0"
$ws.Range("K15").NumberFormat = "@"
$ws.Range("K15").Value = "***.*"
$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = "***.*"
$ws.Range("M15").NumberFormat = "@"
$ws.Range("M15").Value = "***.*"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "0"
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "***.*"
$ws.Range("M16").NumberFormat = "@"
$ws.Range("M16").Value = "***.*"
$ws.Range("M17").NumberFormat = "@"
$ws.Range("M17").Value = "***.*"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "0"
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value = "***.*"
$ws.Range("I20").NumberFormat = "@"
$ws.Range("I20").Value = "0"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "0"
$ws.Range("K20").NumberFormat = "@"
$ws.Range("K20").Value = "***.*"
$ws.Range("L20").NumberFormat = "@"
$ws.Range("L20").Value = "***.*"
$ws.Range("M20").NumberFormat = "@"
$ws.Range("M20").Value = "***.*"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value = "0"
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value = "0"
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = "***.*"
$ws.Range("L22").NumberFormat = "@"
$ws.Range("L22").Value = "***.*"
$ws.Range("M22").NumberFormat = "@"
$ws.Range("M22").Value = "***.*"
$ws.Range("I23").NumberFormat = "@"
$ws.Range("I23").Value = "0"
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = "0"
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = "***.*"
$ws.Range("L23").NumberFormat = "@"
$ws.Range("L23").Value = "***.*"
$ws.Range("M23").NumberFormat = "@"
$ws.Range("M23").Value = "***.*"
$ws.Range("I25").NumberFormat = "@"
$ws.Range("I25").Value = "0"
$ws.Range("J25").NumberFormat = "@"
$ws.Range("J25").Value = "0"
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = "***.*"
$ws.Range("L25").NumberFormat = "@"
$ws.Range("L25").Value = "***.*"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("I26").NumberFormat = "@"
$ws.Range("I26").Value = "0"
$ws.Range("J26").NumberFormat = "@"
$ws.Range("J26").Value = "0"
$ws.Range("K26").NumberFormat = "@"
$ws.Range("K26").Value = "***.*"
$ws.Range("L26").NumberFormat = "@"
$ws.Range("L26").Value = "***.*"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("I27").NumberFormat = "@"
$ws.Range("I27").Value = "0"
$ws.Range("L27").NumberFormat = "@"
$ws.Range("L27").Value = "***.*"
$ws.Range("I28").NumberFormat = "@"
$ws.Range("I28").Value = "0"
$ws.Range("J28").NumberFormat = "@"
$ws.Range("J28").Value = "0"
$ws.Range("K28").NumberFormat = "@"
$ws.Range("K28").Value = "***.*"
$ws.Range("L28").NumberFormat = "@"
$ws.Range("L28").Value = "***.*"
$ws.Range("M28").NumberFormat = "@"
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").NumberFormat = "@"
$ws.Range("N28").Value = "***.*"
$ws.Range("I29").NumberFormat = "@"
$ws.Range("I29").Value = "0"
$ws.Range("J29").NumberFormat = "@"
$ws.Range("J29").Value = "0"
$ws.Range("K29").NumberFormat = "@"
$ws.Range("K29").Value = "***.*"
$ws.Range("L29").NumberFormat = "@"
$ws.Range("L29").Value = "***.*"
$ws.Range("M29").NumberFormat = "@"
$ws.Range("M29").Value = "***.*"
$ws.Range("N29").NumberFormat = "@"
$ws.Range("N29").Value = "***.*"
$ws.Range("I30").NumberFormat = "@"
$ws.Range("I30").Value = "0"
$ws.Range("J30").NumberFormat = "@"
$ws.Range("J30").Value = "0"
$ws.Range("K30").NumberFormat = "@"
$ws.Range("K30").Value = "***.*"
$ws.Range("L30").NumberFormat = "@"
$ws.Range("L30").Value = "***.*"

$ws.Range("A14").Copy()
$ws.Range("I14,J14,K14,L14,M14,N14,D15,E15,I15,J15,K15,L15,M15,I16,L16,M16,M17,I18,L19,I20,J20,K20,L20,M20,C22,I22,J22,K22,L22,M22,I23,J23,K23,L23,M23,I25,J25,K25,L25,D26,E26,I26,J26,K26,L26,C27,I27,L27,I28,J28,K28,L28,M28,N28,I29,J29,K29,L29,M29,N29,I30,J30,K30,L30").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Output "Edit complete"